$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "K_TYPEAREA"
$ws.Range("B13").Value = "Art der Fläche"
$ws.Range("C13").Value = "Type of area"

$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
